# Generate Report for Handoff
# - Flip the localization status from "In Translation" to "Ready for handoff"
#   on the Overview sheet (zh-cn/de-de mirror columns) and on each language
#   sheet's Status column.
# - Bump the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps to reflect the new handoff.
# - Widen the Status-ish columns so the longer "Ready for handoff" text
#   still fits (mirrors Excel's auto-fit after the content grew).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
# Overview's "Latest HO Xliff Generate Date" and de-de's "Latest Handoff
# Datetime" shared the same original value/string and move together.
$wsOverview.Range("G2").Value = "2016-08-19 00:38:07"
$wsDeDe.Range("H2").Value = "2016-08-19 00:38:07"
# zh-cn's "Latest Handoff Datetime" moves independently.
$wsZhCn.Range("H2").Value = "2016-08-19 00:37:59"

# --- Column widths: widen to fit "Ready for handoff" ---
$newColWidth = 16.333333333333332
$wsOverview.Range("E1").ColumnWidth = $newColWidth
$wsOverview.Range("F1").ColumnWidth = $newColWidth
$wsZhCn.Range("C1").ColumnWidth = $newColWidth
$wsDeDe.Range("C1").ColumnWidth = $newColWidth
